$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF -----------------------
# Copy the style of an existing header cell (AC1, style index "1":
# bold font, thin border, centered/top aligned) onto the three new
# header cells before setting their text, so they match the rest of
# row 1's formatting.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-45): season record repeated on every row -----------
$ws.Range("AD2:AD45").Value = 72
$ws.Range("AE2:AE45").Value = 89
$ws.Range("AF2:AF45").Value = 0
